# Re-process the data with the newly curated dimensions.
# Column C = "municipio-nombre" now becomes a dimension (refArea/dim/URI-Municipio)
# instead of a measure, and column D = "diputados" now becomes a measure
# (iaest-measure:diputados / medida / xsd:int) instead of a dimension, and no
# longer has its own mapping workbook reference in row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = "sdmx-dimension:refArea"
$ws.Cells.Item(2,4).Value = "iaest-measure:diputados"

$ws.Cells.Item(3,3).Value = "dim"
$ws.Cells.Item(3,4).Value = "medida"

$ws.Cells.Item(4,3).Value = "URI-Municipio"
$ws.Cells.Item(4,4).Value = "xsd:int"

$ws.Cells.Item(5,4).Clear()
